$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.012.43'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -3.35%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.716.07'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -3.01%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '307.78'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -6.36%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4770'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +4.80%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3470'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.70%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '42.05'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.11%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07234'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -2.07%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.042'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -4.98%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.04%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.80'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -4.45%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.834'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -2.95%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.715.06'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -3.40%  '
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -4.97%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '86.40'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -6.77%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001037'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -2.24%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.07%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.000'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.03%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.48'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -2.72%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.616'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -2.82%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '27.064.28'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -3.23%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.73'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -4.36%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.57%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '19.97'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.96%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '151.06'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -5.84%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.911.94'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -3.29%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.094'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -3.79%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '120.76'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -2.74%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.029'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -4.94%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09146'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.22%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.603'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.60%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.317'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -5.18%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +6.38%  '
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.05877'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -4.08%  '
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02175'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -4.72%  '
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = 'Aptos'
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '10.97'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -7.35%  '
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = 'Algorand'
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.1998'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -4.52%  '
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6016'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -3.93%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.734'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -3.83%  '
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = 'Frax'
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.000'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.07%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.082'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -8.37%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.435'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -4.78%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '12.76'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -3.72%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.565'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -4.50%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5593'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -4.48%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '119.31'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -2.85%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.828'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -5.73%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.111'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.90%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06649'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -2.47%  '
